# Update the "展览" (Exhibition) and "全部类型" (All types) sheets with
# refreshed attendee/view counts in column F for the affected rows.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8357
    3  = 7827
    4  = 128
    8  = 129
    9  = 122
    10 = 169
    13 = 130
    14 = 1358
    16 = 53
    19 = 126
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
